$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 173, pushing rows 173:190 down to
# 174:191 (Excel copies formatting, e.g. the date style on column D, from
# the row above automatically).
$ws.Rows.Item(173).Insert()

# Populate the newly inserted row 173 with the new record.
$ws.Cells.Item(173, 1).Value  = 3
$ws.Cells.Item(173, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(173, 3).Value  = "Coquimbo"
$ws.Cells.Item(173, 4).Value  = 44449
$ws.Cells.Item(173, 5).Value  = 5
$ws.Cells.Item(173, 6).Value  = 100112043
$ws.Cells.Item(173, 7).Value  = "Pepino ensalada"
$ws.Cells.Item(173, 8).Value  = "Sin especificar"
$ws.Cells.Item(173, 9).Value  = "Primera"
$ws.Cells.Item(173, 10).Value = 95
$ws.Cells.Item(173, 11).Value = 14000
$ws.Cells.Item(173, 12).Value = 15000
$ws.Cells.Item(173, 13).Value = 14526
$ws.Cells.Item(173, 14).Value = "$/caja 70 unidades"
$ws.Cells.Item(173, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(173, 16).Value = 208
$ws.Cells.Item(173, 17).Value = 70
$ws.Cells.Item(173, 18).Value = "Hortaliza"
